$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = "Riquelme"
$ws.Range("D7").Value = "Juan Roman"
$ws.Range("E7").Value = "1000000"
$ws.Range("F7").Value = "1010000001"
$ws.Range("G7").Value = "boca@juniors.net"
$ws.Range("H7").Value = "Secundario"
$ws.Range("I7").Value = "NO"
